$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Train Results")
$ws2 = $wb.Worksheets.Item("Test Results")

$sheet1Data = @{
  2 = @(0,40,4,0,28,4,4,20,1.13,1.136703252792358)
  3 = @(4,0,8,24,4,16,44,0,1.22,1.203043937683105)
  4 = @(4,16,0,20,4,12,44,0,1.2,1.197904706001282)
  5 = @(0,24,4,0,8,4,36,24,1.09,1.086240410804749)
  6 = @(4,0,8,4,24,12,44,4,1.13,1.149748921394348)
  7 = @(0,16,8,4,16,20,36,0,1.15,1.157410502433777)
  8 = @(4,20,4,4,16,0,52,0,1.14,1.132760643959045)
  9 = @(4,0,12,4,4,16,52,8,1.16,1.15317702293396)
  10 = @(4,12,8,0,8,8,56.00000000000001,4,1.09,1.114409327507019)
  11 = @(4,4,4,8,28,16,32,4,1.2,1.221900701522827)
  12 = @(4,0,8,4,24,12,44,4,1.17,1.149748921394348)
  13 = @(4,8,0,12,20,8,48,0,1.17,1.125646710395813)
  14 = @(0,20,8,4,0,20,44,4,1.11,1.131448745727539)
  15 = @(0,20,8,4,0,20,44,4,1.13,1.131448745727539)
  16 = @(4,12,0,0,16,8,52,8,1.07,1.095731854438782)
  17 = @(4,12,8,4,16,8,48,0,1.18,1.163147568702698)
  18 = @(0,20,4,0,4,4,48,20,1.15,1.13364851474762)
  19 = @(4,16,0,4,12,8,52,4,1.17,1.159852147102356)
  20 = @(4,8,4,4,16,12,48,4,1.07,1.09487247467041)
  21 = @(4,12,0,12,20,8,44,0,1.09,1.098811030387878)
  22 = @(4,12,0,12,20,8,44,0,1.13,1.098811030387878)
  23 = @(0,12,8,4,16,20,36,4,1.09,1.113114714622498)
  24 = @(4,12,4,4,20,16,32,8,1.17,1.176286697387695)
  25 = @(4,8,0,0,4,8,52,24,1.16,1.161266326904297)
  26 = @(4,0,4,0,20,8,52,12,1.07,1.080740213394165)
  27 = @(4,8,12,4,4,24,39.99999999999999,4,1.14,1.142421364784241)
  28 = @(4,0,12,16,4,12,52,0,1.16,1.158882975578308)
  29 = @(0,16,8,4,16,20,36,0,1.17,1.157410502433777)
  30 = @(4,12,8,0,8,8,56.00000000000001,4,1.12,1.114409327507019)
  31 = @(4,0,4,4,16,12,56.00000000000001,4,1.14,1.139896988868713)
  32 = @(0,4,8,4,16,20,44,4,1.17,1.13602888584137)
  33 = @(0,12,4,0,4,16,39.99999999999999,24,1.17,1.166258335113525)
  34 = @(0,0,8,4,16,12,52,8,1.1,1.120237469673157)
  35 = @(0,16,0,4,20,20,28,12,1.13,1.161486983299255)
  36 = @(4,16,8,0,12,0,48,12,1.08,1.076404929161072)
  37 = @(0,0,8,4,16,12,52,8,1.14,1.120237469673157)
  38 = @(4,12,4,4,12,16,36,12,1.21,1.197951436042786)
  39 = @(0,24,8,0,12,16,32,8,1.15,1.205635070800781)
  40 = @(4,32,8,4,4,20,24,4,1.16,1.161118626594543)
  41 = @(0,24,8,0,12,16,32,8,1.25,1.205635070800781)
  42 = @(4,40,0,0,12,4,36,4,1.18,1.178856253623962)
  43 = @(0,16,0,4,20,20,28,12,1.17,1.161486983299255)
  44 = @(4,20,4,4,16,0,52,0,1.13,1.132760643959045)
  45 = @(4,12,0,0,16,8,52,8,1.13,1.095731854438782)
}

$sheet2Data = @{
  2 = @(0,16,8,4,16,20,36,0,1.17,1.157410502433777)
  3 = @(4,0,8,4,24,12,44,4,1.17,1.149748921394348)
  4 = @(4,12,8,0,8,8,56.00000000000001,4,1.09,1.114409327507019)
  5 = @(0,16,8,4,16,20,36,0,1.15,1.157410502433777)
  6 = @(4,12,0,12,20,8,44,0,1.09,1.098811030387878)
  7 = @(4,12,4,4,12,16,36,12,1.21,1.197951436042786)
  8 = @(0,24,8,0,12,16,32,8,1.25,1.205635070800781)
  9 = @(4,0,4,0,4,4,60,24,1.14,1.129131078720093)
  10 = @(4,0,12,4,4,16,52,8,1.16,1.15317702293396)
  11 = @(0,24,4,0,8,4,36,24,1.09,1.086240291595459)
  12 = @(4,16,8,0,12,0,48,12,1.08,1.076404929161072)
  13 = @(0,12,8,4,16,20,36,4,1.09,1.113114714622498)
}

foreach ($r in $sheet1Data.Keys) {
  $rowVals = $sheet1Data[$r]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws1.Cells.Item($r, $c + 1).Value = $rowVals[$c]
  }
}

foreach ($r in $sheet2Data.Keys) {
  $rowVals = $sheet2Data[$r]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws2.Cells.Item($r, $c + 1).Value = $rowVals[$c]
  }
}
